$wb = $excel.ActiveWorkbook

# This script applies updated market-board price/profit values scraped
# by the scheduled runner. Values are plain data (no formulas in this
# workbook), so each touched cell is written directly. A couple of cells
# are removed (ClearContents) or newly introduced (Value =) to mirror the
# exact shape of the upstream diff.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 505.25
$ws.Range("I53").Value = 195.375
$ws.Range("J53").Value = 1125
$ws.Range("K53").Value = 195.375
$ws.Range("L53").Value = 1125
$ws.Range("M53").Value = 441.625
$ws.Range("N53").Value = -2399
$ws.Range("H54").Value = 6482
$ws.Range("I54").Value = 6482
$ws.Range("K54").Value = 6482
$ws.Range("M54").Value = -5996
$ws.Range("H135").Value = 481.33334
$ws.Range("I135").Value = 481.33334
$ws.Range("K135").Value = 4332.00006
$ws.Range("M135").Value = -1797.00006
$ws.Range("H138").Value = 4743.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1278.5294
$ws.Range("I2").Value = 612.4545000000001
$ws.Range("K2").Value = 612.4545000000001
$ws.Range("M2").Value = -499.4545000000001
$ws.Range("H13").Value = 19800.6
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 19800.6
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 19800.6
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -20088.6
$ws.Range("H32").Value = 36181.29
$ws.Range("I32").Value = 37593.57
$ws.Range("K32").Value = 37593.57
$ws.Range("M32").Value = -37306.57
$ws.Range("H61").Value = 2815.4119
$ws.Range("I61").Value = 2633.0715
$ws.Range("J61").Value = 3666.3333
$ws.Range("K61").Value = 2633.0715
$ws.Range("L61").Value = 3666.3333
$ws.Range("M61").Value = -2421.0715
$ws.Range("N61").Value = -4090.3333
$ws.Range("H88").Value = 2068.7222
$ws.Range("I88").Value = 1425.5714
$ws.Range("K88").Value = 1425.5714
$ws.Range("M88").Value = -1019.5714
$ws.Range("H91").Value = 2068.7222
$ws.Range("I91").Value = 1425.5714
$ws.Range("K91").Value = 1425.5714
$ws.Range("M91").Value = -21.57140000000004
$ws.Range("H116").Value = 1278.5294
$ws.Range("I116").Value = 612.4545000000001
$ws.Range("K116").Value = 612.4545000000001
$ws.Range("M116").Value = 1681.5455
$ws.Range("H122").Value = 3196.6667
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3295
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9885
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -14785
$ws.Range("H136").Value = 2815.4119
$ws.Range("I136").Value = 2633.0715
$ws.Range("J136").Value = 3666.3333
$ws.Range("K136").Value = 7899.2145
$ws.Range("L136").Value = 10998.9999
$ws.Range("M136").Value = -5349.2145
$ws.Range("N136").Value = -16098.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1278.5294
$ws.Range("I3").Value = 612.4545000000001
$ws.Range("K3").Value = 612.4545000000001
$ws.Range("M3").Value = -498.4545000000001
$ws.Range("H5").Value = 2531.8333
$ws.Range("I5").Value = 3672.75
$ws.Range("K5").Value = 3672.75
$ws.Range("M5").Value = -3559.75
$ws.Range("H13").Value = 59999
$ws.Range("J13").Value = 59999
$ws.Range("L13").Value = 59999
$ws.Range("N13").Value = -60335
$ws.Range("H99").Value = 49045.137
$ws.Range("I99").Value = 55212.79
$ws.Range("K99").Value = 55212.79
$ws.Range("M99").Value = -53714.79
$ws.Range("H105").Value = 3429.7896
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 5000
$ws.Range("N105").Value = -8494
$ws.Range("H134").Value = 3153.2856
$ws.Range("I134").Value = 2119.5557
$ws.Range("K134").Value = 6358.6671
$ws.Range("M134").Value = -3823.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 69999
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H31").Value = 1366.3334
$ws.Range("I31").Value = 1366.3334
$ws.Range("K31").Value = 1366.3334
$ws.Range("M31").Value = -1071.3334
$ws.Range("H34").Value = 1366.3334
$ws.Range("I34").Value = 1366.3334
$ws.Range("K34").Value = 1366.3334
$ws.Range("M34").Value = -1164.3334
$ws.Range("H75").Value = 69460
$ws.Range("J75").Value = 69460
$ws.Range("L75").Value = 69460
$ws.Range("N75").Value = -71456
$ws.Range("H78").Value = 69460
$ws.Range("J78").Value = 69460
$ws.Range("L78").Value = 208380
$ws.Range("N78").Value = -218364
$ws.Range("H99").Value = 4636.3335
$ws.Range("I99").Value = 4554.8
$ws.Range("K99").Value = 4554.8
$ws.Range("M99").Value = -3056.8
$ws.Range("H126").Value = 4636.3335
$ws.Range("I126").Value = 4554.8
$ws.Range("K126").Value = 13664.4
$ws.Range("M126").Value = -11194.4
$ws.Range("H132").Value = 1536.5385
$ws.Range("I132").Value = 1372.125
$ws.Range("K132").Value = 4116.375
$ws.Range("M132").Value = -1586.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 673.25
$ws.Range("J25").Value = 673.25
$ws.Range("L25").Value = 2019.75
$ws.Range("N25").Value = -2357.75
$ws.Range("H30").Value = 673.25
$ws.Range("J30").Value = 673.25
$ws.Range("L30").Value = 2019.75
$ws.Range("N30").Value = -2223.75
$ws.Range("H92").Value = 298.57144
$ws.Range("I92").Value = 248.22223
$ws.Range("J92").Value = 389.2
$ws.Range("K92").Value = 744.66669
$ws.Range("L92").Value = 1167.6
$ws.Range("M92").Value = 503.33331
$ws.Range("N92").Value = -3663.6
$ws.Range("H108").Value = 3218.5
$ws.Range("I108").Value = 3218.5
$ws.Range("K108").Value = 9655.5
$ws.Range("M108").Value = -6775.5
$ws.Range("H140").Value = 2621.8572
$ws.Range("I140").Value = 2436.3845
$ws.Range("J140").Value = 5033
$ws.Range("K140").Value = 7309.1535
$ws.Range("L140").Value = 15099
$ws.Range("M140").Value = -2129.1535
$ws.Range("N140").Value = -25459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 95992.336
$ws.Range("J45").Value = 95992.336
$ws.Range("L45").Value = 95992.336
$ws.Range("N45").Value = -97110.336
$ws.Range("H88").Value = 80000
$ws.Range("J88").Value = 80000
$ws.Range("L88").Value = 80000
$ws.Range("N88").Value = -80902
$ws.Range("H91").Value = 80000
$ws.Range("J91").Value = 80000
$ws.Range("L91").Value = 80000
$ws.Range("N91").Value = -83120
$ws.Range("H102").Value = 3498.9
$ws.Range("I102").Value = 2519.3
$ws.Range("J102").Value = 4478.5
$ws.Range("K102").Value = 2519.3
$ws.Range("L102").Value = 4478.5
$ws.Range("M102").Value = -897.3000000000002
$ws.Range("N102").Value = -7722.5
$ws.Range("H106").Value = 24280
$ws.Range("J106").Value = 24280
$ws.Range("L106").Value = 24280
$ws.Range("N106").Value = -26804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2516.3076
$ws.Range("I7").Value = 1892.6666
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 1892.6666
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -1780.6666
$ws.Range("N7").Value = -10224
$ws.Range("H22").Value = 56728.25
$ws.Range("J22").Value = 1560.6154
$ws.Range("L22").Value = 1560.6154
$ws.Range("N22").Value = -2150.6154
$ws.Range("H27").Value = 56728.25
$ws.Range("J27").Value = 1560.6154
$ws.Range("L27").Value = 1560.6154
$ws.Range("N27").Value = -1774.6154
$ws.Range("H61").Value = 8249.5
$ws.Range("I61").Value = 7999
$ws.Range("K61").Value = 7999
$ws.Range("M61").Value = -7797
$ws.Range("H68").Value = 4232.3335
$ws.Range("J68").Value = 5879.4
$ws.Range("L68").Value = 5879.4
$ws.Range("N68").Value = -7377.4
$ws.Range("H71").Value = 4232.3335
$ws.Range("J71").Value = 5879.4
$ws.Range("L71").Value = 29397
$ws.Range("N71").Value = -36885
$ws.Range("H113").Value = 8249.5
$ws.Range("I113").Value = 7999
$ws.Range("K113").Value = 7999
$ws.Range("M113").Value = -5829
$ws.Range("H126").Value = 2516.3076
$ws.Range("I126").Value = 1892.6666
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 5677.9998
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -3207.9998
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 45750.55
$ws.Range("I132").Value = 59608.19
$ws.Range("K132").Value = 178824.57
$ws.Range("M132").Value = -176294.57
$ws.Range("H133").Value = 92461
$ws.Range("J133").Value = 92461
$ws.Range("L133").Value = 92461
$ws.Range("N133").Value = -97521

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2990
$ws.Range("I3").Value = 2990
$ws.Range("K3").Value = 2990
$ws.Range("M3").Value = -2876
$ws.Range("H10").Value = 6000
$ws.Range("I10").Value = 6000
$ws.Range("K10").Value = 6000
$ws.Range("M10").Value = -5831
$ws.Range("H45").Value = 59503.445
$ws.Range("J45").Value = 56504.57
$ws.Range("L45").Value = 56504.57
$ws.Range("N45").Value = -57486.57
$ws.Range("H100").Value = 1271.1428
$ws.Range("I100").Value = 774.5
$ws.Range("J100").Value = 1933.3334
$ws.Range("K100").Value = 1549
$ws.Range("L100").Value = 3866.6668
$ws.Range("M100").Value = -1008
$ws.Range("N100").Value = -4948.6668
$ws.Range("H132").Value = 84906.16
$ws.Range("I132").Value = 87815.086
$ws.Range("K132").Value = 263445.258
$ws.Range("M132").Value = -260915.258

